# Weekly update: insert a new price record for "Cebolla" (Macroferia Regional
# de Talca) at row 656, shifting the existing rows 656-701 down to 657-702.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 656 (shifts rows 656:701 -> 657:702,
# and extends the used range to row 702).
$ws.Rows.Item(656).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(656, 1).Value = 5
$ws.Cells.Item(656, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(656, 3).Value = "Maule"
$ws.Cells.Item(656, 4).Value = 44826
$ws.Cells.Item(656, 5).Value = 7
$ws.Cells.Item(656, 6).Value = 100112004
$ws.Cells.Item(656, 7).Value = "Cebolla"
$ws.Cells.Item(656, 8).Value = "Sin especificar"
$ws.Cells.Item(656, 9).Value = "1a (guarda)"
$ws.Cells.Item(656, 10).Value = 1500
$ws.Cells.Item(656, 11).Value = 13000
$ws.Cells.Item(656, 12).Value = 13000
$ws.Cells.Item(656, 13).Value = 13000
$ws.Cells.Item(656, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(656, 15).Value = "Región del Maule"
$ws.Cells.Item(656, 16).Value = 520
$ws.Cells.Item(656, 17).Value = 25
$ws.Cells.Item(656, 18).Value = "Hortaliza"
